$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 3025
$ws1.Range("F13").Value = 326
$ws1.Range("F22").Value = 6902
$ws1.Range("F23").Value = 6902
$ws1.Range("F27").Value = 1272
$ws1.Range("F38").Value = 6112
$ws1.Range("F49").Value = 360

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3025
$ws4.Range("F13").Value = 326
$ws4.Range("F22").Value = 6902
$ws4.Range("F23").Value = 6902
$ws4.Range("F27").Value = 1272
$ws4.Range("F39").Value = 6112
$ws4.Range("F48").Value = 360
